# ValueSet-eclaire-study-phase-vs.xlsx correction:
#   - Remove the "Include from ResearchStudyPha" include sheet (the old
#     http://terminology.hl7.org/CodeSystem/research-study-phase system
#     is no longer included), leaving only "Metadata" and
#     "Include from Définition des t".
#   - Bump the Metadata "Date" value to reflect the re-generation time.

$wb = $excel.ActiveWorkbook

$obsolete = $wb.Worksheets.Item("Include from ResearchStudyPha")
$obsolete.Delete()

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-09-01T17:43:23+00:00"
